# The workbook "A2_Structure_Lists.xlsx" has a single sheet, "Lists", which
# is a flat table of (Year, Tech, Timeslice, Fuel, Emission, MOO, Region)
# rows used to configure the model. The commit adapts the sheet "to the
# last RD model": the Region code in row 2 (column G) is updated from the
# old "CR" code to the new "RD" code.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")

$ws.Range("G2").Value = "RD"
